$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 9 text is rewritten with the final requirement wording
$ws.Range("B9").Value = "La aplicación permitira fijar un camino habilitado para el jugador"

# New requerimiento deseado (row4, col C) - sits right under the "Requerimientos deseados" header
$ws.Range("C4").Value = "La aplicación podría permitir mostrar un ranking con todos los puntajes"

# New row 10 with another functional requirement
$ws.Range("B10").Value = "La aplicaicon permitira habilitar un punto de llegada para el jugador"

# Column C was manually widened (no longer auto "best fit")
$ws.Columns("C").ColumnWidth = 63

# Selection moved to the next empty row after the edits
[void]$ws.Range("B11").Select()
